$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = "Ejemplo de cosotos unitarios"
$ws.Range("A2").Select()
